# Commit: "complate!!-> scrapping whole ipl"
#
# Changes applied to the "Virat Singh" batting sheet:
#   1. Rename the sheet from "Sheet1" to "Virat Singh".
#   2. Insert a new "matchNo" column at the very left (column A), shifting
#      every existing column one place to the right.
#   3. Fill in the matchNo value ("20th") for the previously scraped match.
#   4. Append a newly scraped match as row 3 ("9th" match vs Mumbai Indians).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab
$ws.Name = "Virat Singh"

# 2. Insert new column A ("matchNo"), pushing teamName..result to B..M
$ws.Columns.Item(1).Insert()

# New header cell
$ws.Range("A1").Value = "matchNo"

# 3. matchNo for the existing (row 2) scraped entry
$ws.Range("A2").Value = "20th"

# 4. Newly scraped match appended as row 3
#    (leading "'" keeps numeric-looking values stored as text, matching
#     the rest of the sheet, which is entirely text-typed)
$ws.Range("A3").Value = "9th"
$ws.Range("B3").Value = "Sunrisers Hyderabad"
$ws.Range("C3").Value = "Virat Singh"
$ws.Range("D3").Value = "c Yadav b Chahar"
$ws.Range("E3").Value = "'11"
$ws.Range("F3").Value = "'12"
$ws.Range("G3").Value = "'1"
$ws.Range("H3").Value = "'0"
$ws.Range("I3").Value = "'91.66"
$ws.Range("J3").Value = "Mumbai Indians"
$ws.Range("K3").Value = "Chennai"
$ws.Range("L3").Value = "April 17"
$ws.Range("M3").Value = "Mumbai won by 13 runs"
